$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 4: "Dados" — dataset size description ("mais de" -> "aproximadamente")
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(2)
$tr4 = $shape4.TextFrame.TextRange
$para4_3 = $tr4.Paragraphs(3)
$para4_3.Runs(3).Text = " utilizado conta com aproximadamente 3 milhões de linhas e pesa mais de 600MB"

# ---------------------------------------------------------------------------
# Slide 15: "Apriori" — add recovery percentage parenthetical
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$shape15 = $s15.Shapes.Item(2)
$tr15 = $shape15.TextFrame.TextRange
$para15_6 = $tr15.Paragraphs(6)
$para15_6.Runs(1).Text = "Isto acontece devido a quantidade de casos nos quais a evolução é “Recuperado” ser muito maior que a de óbitos (>98.5% de recuperados)."

# ---------------------------------------------------------------------------
# Slide 24: "Resultados obtidos" — expand first bullet, rewrite second bullet
# ---------------------------------------------------------------------------
$s24 = $p.Slides.Item(24)
$shape24 = $s24.Shapes.Item(2)
$tr24 = $shape24.TextFrame.TextRange

$para24_1 = $tr24.Paragraphs(1)
$para24_1.Runs(1).Text = "Os resultados obtidos com os métodos de mineração de dados mostraram padrões nos dados que não são triviais, como por exemplo que alguns sintomas prevalecem em comparação a outros quando se trata de pacientes que foram a óbito."

# Paragraph 2 originally has two runs ("Isto pode ser útil " + "em casos ");
# remove it entirely and replace with a single fresh paragraph inheriting
# paragraph 1's (clean) run formatting.
$para24_2 = $tr24.Paragraphs(2)
$para24_2.Delete()
$para24_1 = $tr24.Paragraphs(1)
[void]$para24_1.InsertAfter([char]13 + "Essas descobertas podem ser valiosas, por exemplo, em situações clínicas, auxiliando na identificação precoce e manejo de pacientes com suspeita de COVID-19.")

# ---------------------------------------------------------------------------
# Slide 25: "Trabalhos futuros" — trailing space + replace "[...]" placeholder
# with two concrete future-work bullets (removing the red highlight color).
# ---------------------------------------------------------------------------
$s25 = $p.Slides.Item(25)
$shape25 = $s25.Shapes.Item(2)
$tr25 = $shape25.TextFrame.TextRange

$para25_1 = $tr25.Paragraphs(1)
$para25_1.Runs(1).Text = "Em trabalhos futuros buscamos explorar ainda mais os dados da doença no estado. "

$para25_1 = $tr25.Paragraphs(1)
$para25_1.InsertAfter([char]13 + "Investigar a evolução dos padrões de sintomas ao longo do tempo." + [char]13 + "Combinar conjuntos de dados da COVID-19 com informações como ocupação de leitos hospitalares e dados de vacinação")

# Remove the old "[...]" (red) paragraph, now in the middle of the list.
$oldPara25 = $tr25.Paragraphs(4)
$oldPara25.Delete()

# Deleting the shape's final paragraph leaves an empty ghost paragraph mark
# behind; a second Delete() on that now-empty trailing paragraph removes it.
$ghost25 = $tr25.Paragraphs(4)
$ghost25.Delete()
